$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''61.869.89'
$ws.Range('E2').Value = '''  +5.35%  '

$ws.Range('D3').Value = '''3.061.74'
$ws.Range('E3').Value = '''  +2.88%  '

$ws.Range('E4').Value = '''  +0.07%  '

$ws.Range('D5').Value = '''577.10'
$ws.Range('E5').Value = '''  +3.01%  '

$ws.Range('D6').Value = '''142.59'
$ws.Range('E6').Value = '''  +4.40%  '

$ws.Range('E7').Value = '''  -0.01%  '

$ws.Range('D8').Value = '''3.048.85'
$ws.Range('E8').Value = '''  +2.72%  '

$ws.Range('D9').Value = '''0.524'
$ws.Range('E9').Value = '''  +1.38%  '

$ws.Range('D10').Value = '''0.139'
$ws.Range('E10').Value = '''  +5.99%  '

$ws.Range('D11').Value = '''5.47'
$ws.Range('E11').Value = '''  +13.04%  '

$ws.Range('D12').Value = '''0.463'
$ws.Range('E12').Value = '''  +1.88%  '

$ws.Range('D13').Value = '''0.0000240'
$ws.Range('E13').Value = '''  +5.87%  '

$ws.Range('D14').Value = '''34.73'
$ws.Range('E14').Value = '''  +4.01%  '

$ws.Range('E15').Value = '''  +0.07%  '

$ws.Range('D16').Value = '''3.565.25'
$ws.Range('E16').Value = '''  +3.01%  '

$ws.Range('E17').Value = '''  +3.39%  '

$ws.Range('D18').Value = '''3.053.78'
$ws.Range('E18').Value = '''  +3.03%  '

$ws.Range('D19').Value = '''61.811.44'
$ws.Range('E19').Value = '''  +5.28%  '

$ws.Range('D20').Value = '''448.68'
$ws.Range('E20').Value = '''  +6.53%  '

$ws.Range('D21').Value = '''13.92'
$ws.Range('E21').Value = '''  +3.21%  '

$ws.Range('D22').Value = '''0.730'
$ws.Range('E22').Value = '''  +2.90%  '

$ws.Range('D23').Value = '''7.28'
$ws.Range('E23').Value = '''  +2.72%  '

$ws.Range('D24').Value = '''13.66'
$ws.Range('E24').Value = '''  +2.57%  '

$ws.Range('D25').Value = '''81.66'
$ws.Range('E25').Value = '''  +1.99%  '

$ws.Range('E26').Value = '''  +0.19%  '

$ws.Range('E27').Value = '''  +6.17%  '

$ws.Range('E28').Value = '''  +0.04%  '

$ws.Range('E29').Value = '''  +4.80%  '

$ws.Range('D30').Value = '''7.99'
$ws.Range('E30').Value = '''  +3.69%  '

$ws.Range('D31').Value = '''6.57'
$ws.Range('E31').Value = '''  +8.50%  '

$ws.Range('D32').Value = '''26.59'
$ws.Range('E32').Value = '''  +4.03%  '

$ws.Range('E33').Value = '''  +7.39%  '

$ws.Range('D34').Value = '''0.0₃0811'
$ws.Range('E34').Value = '''  +7.05%  '

$ws.Range('E35').Value = '''  +2.81%  '

$ws.Range('D36').Value = '''6.07'
$ws.Range('E36').Value = '''  +6.14%  '

$ws.Range('D37').Value = '''2.18'
$ws.Range('E37').Value = '''  +5.68%  '

$ws.Range('D38').Value = '''50.18'
$ws.Range('E38').Value = '''  +2.85%  '

$ws.Range('D39').Value = '''2.96'
$ws.Range('E39').Value = '''  +7.72%  '

$ws.Range('E40').Value = '''  +2.29%  '

$ws.Range('D41').Value = '''413.70'
$ws.Range('E41').Value = '''  +4.26%  '

$ws.Range('D42').Value = '''0.0367'
$ws.Range('E42').Value = '''  +6.00%  '

$ws.Range('D43').Value = '''2.766.03'
$ws.Range('E43').Value = '''  +0.86%  '

$ws.Range('E44').Value = '''  +0.69%  '

$ws.Range('D45').Value = '''0.264'
$ws.Range('E45').Value = '''  +9.17%  '

$ws.Range('D46').Value = '''36.82'
$ws.Range('E46').Value = '''  +14.90%  '

$ws.Range('E47').Value = '''  -0.03%  '

$ws.Range('E48').Value = '''  +4.13%  '

$ws.Range('D49').Value = '''122.99'
$ws.Range('E49').Value = '''  -1.52%  '

$ws.Range('E50').Value = '''  +1.69%  '

$ws.Range('D51').Value = '''24.01'
$ws.Range('E51').Value = '''  +3.81%  '
